# "more picks, updated depth charts"
#
# 1. Depth chart update: balco's Ketel Marte (row 6) moves from OF to 2B.
# 2. Eight new draft picks appended as rows 413-420 on the "draftpicks" sheet.
# 3. Viewport/selection nudged down to the newly-added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("draftpicks")

# --- updated depth chart --------------------------------------------------
$ws.Range("D6").Value = "2B"

# --- new picks -------------------------------------------------------------
# NB: shared-string indices are assigned in first-write order, and the
# canonical file has "Seth Lugo" (row 420) land before "Wilmer Flores"
# (row 413) in the shared string table, so B420 is written first.
$ws.Range("B420").Value = "Seth Lugo"
$ws.Range("B413").Value = "Wilmer Flores"
$ws.Range("B414").Value = "Clint Frazier"
$ws.Range("B415").Value = "Adrian Houser"
$ws.Range("B416").Value = "Shed Long"
$ws.Range("B417").Value = "Brandon Belt"
$ws.Range("B418").Value = "Jasson Dominguez"
$ws.Range("B419").Value = "Nick Anderson"

$ws.Range("A413").Value = "marmaduke"
$ws.Range("C413").Value = 1
$ws.Range("D413").Value = "MI"
$ws.Range("E413").Value = 43871

$ws.Range("A414").Value = "balco"
$ws.Range("C414").Value = 3
$ws.Range("D414").Value = "OF"
$ws.Range("E414").Value = 43871

$ws.Range("A415").Value = "sturgeon"
$ws.Range("C415").Value = 5
$ws.Range("D415").Value = "P"
$ws.Range("E415").Value = 43871

$ws.Range("A416").Value = "drjames"
$ws.Range("C416").Value = 2
$ws.Range("D416").Value = "OF"
$ws.Range("E416").Value = 43871

$ws.Range("A417").Value = "drjames"
$ws.Range("C417").Value = 2
$ws.Range("D417").Value = "CI"
$ws.Range("E417").Value = 43871

$ws.Range("A418").Value = "sturgeon"
$ws.Range("C418").Value = 1
$ws.Range("D418").Value = "DH"
$ws.Range("E418").Value = 43871

$ws.Range("A419").Value = "marmaduke"
$ws.Range("C419").Value = 13
$ws.Range("D419").Value = "P"
$ws.Range("E419").Value = 43871

$ws.Range("A420").Value = "bears"
$ws.Range("C420").Value = 3
$ws.Range("D420").Value = "P"
$ws.Range("E420").Value = 43871

# --- viewport: scroll down toward the newly added rows and move the
# active selection to reflect where editing left off -----------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 408
$win.ScrollColumn = 1
$ws.Range("B422").Select()
